# Update "想去人数" (interest count) values that were refreshed in the
# upstream data source. Each sheet below is updated in-place by cell
# reference with the new value taken from the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2  = 37570
$ws.Range("F10").Value2 = 92
$ws.Range("F12").Value2 = 534
$ws.Range("F13").Value2 = 35
$ws.Range("F23").Value2 = 2516
$ws.Range("F24").Value2 = 999
$ws.Range("F28").Value2 = 45
$ws.Range("F29").Value2 = 766
$ws.Range("F31").Value2 = 1153

# --- Sheet "演出" -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value2  = 385
$ws.Range("F5").Value2  = 2
$ws.Range("F9").Value2  = 142
$ws.Range("F10").Value2 = 11

# --- Sheet "全部类型" --------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value2  = 37570
$ws.Range("F11").Value2 = 385
$ws.Range("F13").Value2 = 2
$ws.Range("F16").Value2 = 92
$ws.Range("F18").Value2 = 534
$ws.Range("F19").Value2 = 35
$ws.Range("F23").Value2 = 142
$ws.Range("F24").Value2 = 11
$ws.Range("F34").Value2 = 2516
$ws.Range("F35").Value2 = 999
$ws.Range("F39").Value2 = 45
$ws.Range("F41").Value2 = 766
$ws.Range("F43").Value2 = 1153
